# Figure_S2.pptx edit:
#   "Updated canadensis to canadense"
#
# 1) Fix the species-name typo " Chl. canadensis L304-6D"" -> " Chl. canadense L304-6D""
#    in the text box on the slide (the run lives inside a nested group).
# 2) The deck's "datetimeFigureOut" footer field was re-cached by PowerPoint at
#    save time (2020-01-10 -> 2020-02-18); that field lives on the slide master
#    and every slide layout, so update it everywhere it appears.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) canadensis -> canadense
# ---------------------------------------------------------------------------
function Find-ShapeWithText($shapes, $needle) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 6) {
            # msoGroup - recurse into it
            $found = Find-ShapeWithText $sh.GroupItems $needle
            if ($found -ne $null) {
                return $found
            }
        } elseif ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -like ("*" + $needle + "*")) {
                return $sh
            }
        }
    }
    return $null
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    $target = Find-ShapeWithText $slide.Shapes "canadensis"
    if ($target -ne $null) {
        $tr = $target.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf("canadensis")
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, "canadensis".Length)
            $sub.Text = "canadense"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Re-cached datetimeFigureOut footer field: 2020-01-10 -> 2020-02-18
#    (slide master + every custom layout)
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes, $newDate) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = -1
        try { $phType = $sh.PlaceholderFormat.Type } catch {}
        if ($phType -eq 16) {
            if ($sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "2020-01-10") {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes "2020-02-18"
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes "2020-02-18"
}
